$wb = $excel.ActiveWorkbook

# --- "AddOpportunity" sheet: update values and add a new "WomenLed" column ---
$wsOpp = $wb.Worksheets.Item("AddOpportunity")

# MarketCap value for the single data row changes from 10.0 to 10000.0
$wsOpp.Range("AA2").Value = "10000.0"

# New header column AD1 = "WomenLed" (bold, like the other header cells)
$wsOpp.Range("AD1").Value = "WomenLed"
$wsOpp.Range("AD1").Font.Bold = $true

# New data value under the WomenLed column
$wsOpp.Range("AD2").Value = "No"

# Contact (Staff) for the row changes from Sam Rogers to Chris Cessna
$wsOpp.Range("N2").Value = "Chris Cessna"

# --- "Users" sheet: move the active selection ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("E15").Select() | Out-Null
